# EmailNoMatch.xlsx - "Added and updated Negative Test Cases for RAD"
#
# The RAD test-data table on Sheet1 is updated:
#  - the "Personal Income Tax" row is removed
#  - the "Sales and Use" row is replaced by a new "New Liability" row that
#    only carries values in the Execute/TaxType columns (Date/Result blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 3 ("Personal Income Tax"); rows below shift up.
$ws.Rows(3).Delete()

# The former "Sales and Use" row is now row 5. Replace it with the new
# "New Liability" row: blank Result/Date, keep Execute = "Y", new TaxType.
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("D5").Value = "New Liability"

# Match the saved selection state.
$ws.Range("A5:B5").Select()
